$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in E10 (855528 -> 854877)
$ws.Range("E10").Value = 854877

# Copy the formatting of row 10 down to new row 11 (keeps A11 bold/bordered/centered
# style consistent with the other "season index" cells, matches column A's style)
$ws.Range("A10:H10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Populate new row 11 values
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "M2_10 Cat 2020"
$ws.Range("C11").Value = 9703
$ws.Range("D11").Value = 10804
$ws.Range("E11").Value = 929613
$ws.Range("F11").Value = 9977
$ws.Range("G11").Value = 10067
$ws.Range("H11").Value = 10176
